# Weekly price-data refresh: insert a new observation row at row 45
# (pushing the existing rows 45-122 down to 46-123), then populate the
# new row with the latest weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 45; existing rows shift down.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new weekly data point.
$ws.Range("A45").Value = 9
$ws.Range("B45").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C45").Value = "Metropolitana"
$ws.Range("D45").Value = 45238
$ws.Range("E45").Value = 13
$ws.Range("F45").Value = 100112029
$ws.Range("G45").Value = "Orégano"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 16
$ws.Range("K45").Value = 21000
$ws.Range("L45").Value = 21000
$ws.Range("M45").Value = 21000
$ws.Range("N45").Value = "$/docena de atados"
$ws.Range("O45").Value = "Región Metropolitana"
$ws.Range("P45").Value = 7000
$ws.Range("Q45").Value = 3
$ws.Range("R45").Value = "Hortaliza"
